# Update TPM-derived values in the LR-pairs sheet (Gnai2-Adcy1) with new TPM-based figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("O2").Value = 0.4600512618675063
$ws.Range("P2").Value = 0.4600512618675063
$ws.Range("Q2").Value = 3.155689493284889
$ws.Range("R2").Value = 28.401205439564
$ws.Range("S2").Value = 0.1558977441612797
$ws.Range("T2").Value = 0.1558977441612796

$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("M3").Value = 0.018749
$ws.Range("N3").Value = 0.05624700000000001
$ws.Range("O3").Value = 0.5399487381324938
$ws.Range("P3").Value = 0.5399487381324937
$ws.Range("Q3").Value = 3.703740650379667
$ws.Range("R3").Value = 33.333665853417
$ws.Range("S3").Value = 0.1829726319973186
$ws.Range("T3").Value = 0.1829726319973186

$ws.Range("I4").Value = 0.1369374790620155
$ws.Range("J4").Value = 0.1369374790620154
$ws.Range("O4").Value = 0.4600512618675063
$ws.Range("P4").Value = 0.4600512618675063
$ws.Range("S4").Value = 0.06299826003943544
$ws.Range("T4").Value = 0.06299826003943543

$ws.Range("I5").Value = 0.1369374790620155
$ws.Range("J5").Value = 0.1369374790620154
$ws.Range("M5").Value = 0.018749
$ws.Range("N5").Value = 0.05624700000000001
$ws.Range("O5").Value = 0.5399487381324938
$ws.Range("P5").Value = 0.5399487381324937
$ws.Range("Q5").Value = 1.496681160247334
$ws.Range("R5").Value = 13.470130442226
$ws.Range("S5").Value = 0.07393921902258005
$ws.Range("T5").Value = 0.07393921902258001

$ws.Range("G6").Value = 148.824417
$ws.Range("H6").Value = 446.473251
$ws.Range("I6").Value = 0.2552967790580629
$ws.Range("J6").Value = 0.2552967790580629
$ws.Range("O6").Value = 0.4600512618675063
$ws.Range("P6").Value = 0.4600512618675063
$ws.Range("Q6").Value = 2.377420453436001
$ws.Range("R6").Value = 21.396784080924
$ws.Range("S6").Value = 0.1174496053563718
$ws.Range("T6").Value = 0.1174496053563718

$ws.Range("G7").Value = 148.824417
$ws.Range("H7").Value = 446.473251
$ws.Range("I7").Value = 0.2552967790580629
$ws.Range("J7").Value = 0.2552967790580629
$ws.Range("M7").Value = 0.018749
$ws.Range("N7").Value = 0.05624700000000001
$ws.Range("O7").Value = 0.5399487381324938
$ws.Range("P7").Value = 0.5399487381324937
$ws.Range("Q7").Value = 2.790308994333
$ws.Range("R7").Value = 25.112780948997
$ws.Range("S7").Value = 0.1378471737016911
$ws.Range("T7").Value = 0.1378471737016911

$ws.Range("G8").Value = 35.426853
$ws.Range("H8").Value = 106.280559
$ws.Range("I8").Value = 0.06077202683121193
$ws.Range("J8").Value = 0.06077202683121192
$ws.Range("O8").Value = 0.4600512618675063
$ws.Range("P8").Value = 0.4600512618675063
$ws.Range("Q8").Value = 0.5659321677240001
$ws.Range("R8").Value = 5.093389509516
$ws.Range("S8").Value = 0.027958247629945
$ws.Range("T8").Value = 0.02795824762994499

$ws.Range("G9").Value = 35.426853
$ws.Range("H9").Value = 106.280559
$ws.Range("I9").Value = 0.06077202683121193
$ws.Range("J9").Value = 0.06077202683121192
$ws.Range("M9").Value = 0.018749
$ws.Range("N9").Value = 0.05624700000000001
$ws.Range("O9").Value = 0.5399487381324938
$ws.Range("P9").Value = 0.5399487381324937
$ws.Range("Q9").Value = 0.6642180668970001
$ws.Range("R9").Value = 5.977962602073
$ws.Range("S9").Value = 0.03281377920126693
$ws.Range("T9").Value = 0.03281377920126692

$ws.Range("G10").Value = 121.3248153333333
$ws.Range("H10").Value = 363.974446
$ws.Range("I10").Value = 0.2081233388901116
$ws.Range("J10").Value = 0.2081233388901115
$ws.Range("O10").Value = 0.4600512618675063
$ws.Range("P10").Value = 0.4600512618675063
$ws.Range("Q10").Value = 1.938123483344889
$ws.Range("R10").Value = 17.443111350104
$ws.Range("S10").Value = 0.09574740468047448
$ws.Range("T10").Value = 0.09574740468047445

$ws.Range("G11").Value = 121.3248153333333
$ws.Range("H11").Value = 363.974446
$ws.Range("I11").Value = 0.2081233388901116
$ws.Range("J11").Value = 0.2081233388901115
$ws.Range("M11").Value = 0.018749
$ws.Range("N11").Value = 0.05624700000000001
$ws.Range("O11").Value = 0.5399487381324938
$ws.Range("P11").Value = 0.5399487381324937
$ws.Range("Q11").Value = 2.274718962684667
$ws.Range("R11").Value = 20.472470664162
$ws.Range("S11").Value = 0.1123759342096371
$ws.Range("T11").Value = 0.1123759342096371

